# Apply the fixes described in the commit message ("fixes #6 and #5")
# Changes:
#  - B4: 5 -> 0
#  - B6: 2.5499999999999998 -> 0.2
#  - B18: (empty) -> -40
# All dependent formulas (D4, D6, D14, E14, D18, D26, E26) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 0
$ws.Range("B6").Value = 0.2
$ws.Range("B18").Value = -40

$excel.CalculateFullRebuild()
$wb.Save()
